$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '247.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '1BNBBNBBestin24h'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.32'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.519'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05617'
$ws.Range('D5').Style = 'Normal'
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.380'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '5GateTokenGT'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.473'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8071'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '7MXTokenMX'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.035'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '8FTXTokenFTT'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1431'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07334'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03114'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'ProBitToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.1257'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '12ProBitTokenPROB'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.02922'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09274'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001661'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'MCDex'
$ws.Range('C17').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.233'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '16MCDexMCB'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.04751'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('B19').Value = 'One'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0005819'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '18OneONEWorstin24h'
$ws.Range('B20').Value = 'TigerCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.006370'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '19TigerCashTCH'
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.005058'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('B22').Value = 'BitKan'
$ws.Range('C22').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.001052'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '21BitKanKAN'
$ws.Range('B23').Value = 'NitroEx'
$ws.Range('C23').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0001503'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '22NitroExNTX'
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.984'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('B25').Value = 'BTSEToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.193'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '24BTSETokenBTSE'
$ws.Range('B26').Value = 'BitpandaEcosystemToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.3268'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '25BitpandaEcosystemTokenBEST'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0003306'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04143'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007092'
$ws.Range('D41').Style = 'Normal'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1039'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002975'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008709'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005647'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000751'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6811'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.01625'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '47BOLOBOLO'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002103'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.01012'
$ws.Range('D50').Style = 'Normal'
